$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 159, shifting existing rows 159-215 down to 160-216
$ws.Rows("159:159").Insert()

# Populate the newly inserted row 159 with the new record's data
$ws.Range("A159").Value2 = 8
$ws.Range("B159").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C159").Value2 = "Coquimbo"
$ws.Range("D159").Value2 = 44524
$ws.Range("E159").Value2 = 4
$ws.Range("F159").Value2 = 100112032
$ws.Range("G159").Value2 = "Zapallo italiano"
$ws.Range("H159").Value2 = "Sin especificar"
$ws.Range("I159").Value2 = "Primera"
$ws.Range("J159").Value2 = 480
$ws.Range("K159").Value2 = 10000
$ws.Range("L159").Value2 = 11000
$ws.Range("M159").Value2 = 10500
$ws.Range("N159").Value2 = '$/caja 70 unidades'
$ws.Range("O159").Value2 = "Provincia de Limarí"
$ws.Range("P159").Value2 = 150
$ws.Range("Q159").Value2 = 70
$ws.Range("R159").Value2 = "Hortaliza"
